$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 new values (previously row 5's data, plus a new B value)
$ws.Range("A3").Value = 112105307
$ws.Range("B3").Value = 89100
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 5754
$ws.Range("F3").Value = "Gultoppig fingersvamp"
$ws.Range("G3").Value = "Ramaria testaceoflava"
$ws.Range("H3").Value = "(Bres.) Corner"
$ws.Range("P3").Value = "Landverktjärnen (Landverktjärnen), Jmt"
$ws.Range("Q3").Value = 446544
$ws.Range("R3").Value = 7032738

# Row 4: only B changes
$ws.Range("B4").Value = 84929

# Row 5 new values (previously row 3's data, plus a new B value)
$ws.Range("A5").Value = 112105682
$ws.Range("B5").Value = 89090
$ws.Range("D5").Value = "VU"
$ws.Range("E5").Value = 5747
$ws.Range("F5").Value = "Läderdoftande fingersvamp"
$ws.Range("G5").Value = "Ramaria safraniolens"
$ws.Range("H5").Value = "Christian"
$ws.Range("P5").Value = "Svensbergsbäcken (Svensbergsbäcken), Jmt"
$ws.Range("Q5").Value = 446627
$ws.Range("R5").Value = 7032919

# Row 6: only B changes
$ws.Range("B6").Value = 90821
